$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Eliminar los periodos de mora anteriores y agregar los nuevos (2202-2209)
# para cada trabajador, manteniendo el Valor Mora especial (34666) solo en
# el periodo mas reciente (2209) de cada uno.

# --- JAIDER ENRIQUE TORRES VILORIA (73434587) ---
$ws.Cells.Item(16, 3).Value = "73434587"
$ws.Cells.Item(16, 4).Value = "JAIDER ENRIQUE TORRES VILORIA"
$ws.Cells.Item(16, 5).Value = "2209"
$ws.Cells.Item(16, 6).Value = 34666

$ws.Cells.Item(17, 3).Value = "73434587"
$ws.Cells.Item(17, 4).Value = "JAIDER ENRIQUE TORRES VILORIA"
$ws.Cells.Item(17, 5).Value = "2208"
$ws.Cells.Item(17, 6).Value = 40000

$ws.Cells.Item(18, 3).Value = "73434587"
$ws.Cells.Item(18, 4).Value = "JAIDER ENRIQUE TORRES VILORIA"
$ws.Cells.Item(18, 5).Value = "2207"
$ws.Cells.Item(18, 6).Value = 40000

$ws.Cells.Item(19, 3).Value = "73434587"
$ws.Cells.Item(19, 4).Value = "JAIDER ENRIQUE TORRES VILORIA"
$ws.Cells.Item(19, 5).Value = "2206"
$ws.Cells.Item(19, 6).Value = 40000

$ws.Cells.Item(20, 3).Value = "73434587"
$ws.Cells.Item(20, 4).Value = "JAIDER ENRIQUE TORRES VILORIA"
$ws.Cells.Item(20, 5).Value = "2205"
$ws.Cells.Item(20, 6).Value = 40000

$ws.Cells.Item(21, 3).Value = "73434587"
$ws.Cells.Item(21, 4).Value = "JAIDER ENRIQUE TORRES VILORIA"
$ws.Cells.Item(21, 5).Value = "2204"
$ws.Cells.Item(21, 6).Value = 40000

$ws.Cells.Item(22, 3).Value = "73434587"
$ws.Cells.Item(22, 4).Value = "JAIDER ENRIQUE TORRES VILORIA"
$ws.Cells.Item(22, 5).Value = "2203"
$ws.Cells.Item(22, 6).Value = 40000

$ws.Cells.Item(23, 3).Value = "73434587"
$ws.Cells.Item(23, 4).Value = "JAIDER ENRIQUE TORRES VILORIA"
$ws.Cells.Item(23, 5).Value = "2202"
$ws.Cells.Item(23, 6).Value = 40000

# --- EDER LUIS TORRES LAMBRAÃ?O (73549649) ---
$ws.Cells.Item(24, 3).Value = "73549649"
$ws.Cells.Item(24, 4).Value = "EDER LUIS TORRES LAMBRAÃ?O"
$ws.Cells.Item(24, 5).Value = "2209"
$ws.Cells.Item(24, 6).Value = 34666

$ws.Cells.Item(25, 3).Value = "73549649"
$ws.Cells.Item(25, 4).Value = "EDER LUIS TORRES LAMBRAÃ?O"
$ws.Cells.Item(25, 5).Value = "2208"
$ws.Cells.Item(25, 6).Value = 40000

$ws.Cells.Item(26, 3).Value = "73549649"
$ws.Cells.Item(26, 4).Value = "EDER LUIS TORRES LAMBRAÃ?O"
$ws.Cells.Item(26, 5).Value = "2207"
$ws.Cells.Item(26, 6).Value = 40000

$ws.Cells.Item(27, 3).Value = "73549649"
$ws.Cells.Item(27, 4).Value = "EDER LUIS TORRES LAMBRAÃ?O"
$ws.Cells.Item(27, 5).Value = "2206"
$ws.Cells.Item(27, 6).Value = 40000

$ws.Cells.Item(28, 3).Value = "73549649"
$ws.Cells.Item(28, 4).Value = "EDER LUIS TORRES LAMBRAÃ?O"
$ws.Cells.Item(28, 5).Value = "2205"
$ws.Cells.Item(28, 6).Value = 40000
